$wb = $excel.ActiveWorkbook

# "展览" sheet: update want-to-go counts for rows 4 and 5
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 38
$ws1.Range("F5").Value = 4941

# "全部类型" sheet: same two events are duplicated at rows 8 and 9
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 38
$ws4.Range("F9").Value = 4941
